$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be written as TEXT (matching the
# original inline-string cell type) without leaving a residual cell style,
# by temporarily switching the cell to Text format, assigning the value, then
# resetting the cell style back to Normal (General).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "31.297.21"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3
$ws.Range("D3").Value = "2.005.12"
$ws.Range("E3").Value = "  +7.11%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
Set-TextValue $ws.Range("D5") "0.7731"
$ws.Range("E5").Value = "  +63.92%  "

# Row 6
Set-TextValue $ws.Range("D6") "258.62"
$ws.Range("E6").Value = "  +5.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3538"
$ws.Range("E8").Value = "  +23.22%  "

# Row 9
Set-TextValue $ws.Range("D9") "28.66"
$ws.Range("E9").Value = "  +31.35%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07056"
$ws.Range("E10").Value = "  +8.69%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.8591"
$ws.Range("E11").Value = "  +18.26%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.08212"
$ws.Range("E12").Value = "  +5.35%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D13") "101.71"
$ws.Range("E13").Value = "  +1.56%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.006.24"
$ws.Range("E14").Value = "  +7.17%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.608"
$ws.Range("E15").Value = "  +8.46%  "

# Row 16
Set-TextValue $ws.Range("D16") "15.66"
$ws.Range("E16").Value = "  +19.67%  "

# Row 17
Set-TextValue $ws.Range("D17") "273.77"
$ws.Range("E17").Value = "  -3.97%  "

# Row 18
$ws.Range("D18").Value = "31.299.73"
$ws.Range("E18").Value = "  +3.03%  "

# Row 19
Set-TextValue $ws.Range("D19") "5.931"
$ws.Range("E19").Value = "  +11.00%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.000007988"
$ws.Range("E20").Value = "  +6.58%  "

# Row 21
$ws.Range("D21").Value = "2.269.11"
$ws.Range("E21").Value = "  +7.40%  "

# Row 22
Set-TextValue $ws.Range("D22") "1.003"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.003"
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
Set-TextValue $ws.Range("D24") "7.166"
$ws.Range("E24").Value = "  +13.48%  "

# Row 25
$ws.Range("E25").Value = "  +11.43%  "

# Row 26
Set-TextValue $ws.Range("D26") "164.86"
$ws.Range("E26").Value = "  +0.88%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.1471"
$ws.Range("E27").Value = "  +52.14%  "

# Row 28
$ws.Range("E28").Value = "  +5.58%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.393"
$ws.Range("E29").Value = "  +26.23%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.621"
$ws.Range("E30").Value = "  +9.00%  "

# Row 31
Set-TextValue $ws.Range("D31") "4.641"
$ws.Range("E31").Value = "  +9.94%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.360"
$ws.Range("E32").Value = "  +2.95%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.435"
$ws.Range("E33").Value = "  +7.05%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.05229"
$ws.Range("E34").Value = "  +8.64%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D35") "1.221"
$ws.Range("E35").Value = "  +8.56%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.7770"
$ws.Range("E36").Value = "  +12.84%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.818"
$ws.Range("E37").Value = "  +3.41%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.01999"
$ws.Range("E38").Value = "  +5.24%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.944"
$ws.Range("E39").Value = "  +3.60%  "

# Row 40
Set-TextValue $ws.Range("D40") "6.718"
$ws.Range("E40").Value = "  +6.94%  "

# Row 41
Set-TextValue $ws.Range("D41") "80.00"
$ws.Range("E41").Value = "  +4.90%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.4722"
$ws.Range("E42").Value = "  +11.93%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.148"
$ws.Range("E43").Value = "  +9.63%  "

# Row 44
Set-TextValue $ws.Range("D44") "107.47"
$ws.Range("E44").Value = "  +6.31%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.8608"
$ws.Range("E45").Value = "  +4.76%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.004"
$ws.Range("E46").Value = "  +0.42%  "

# Row 47
$ws.Range("E47").Value = "  +11.12%  "

# Row 48
Set-TextValue $ws.Range("D48") "9.911"
$ws.Range("E48").Value = "  +1.31%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.4342"
$ws.Range("E49").Value = "  +11.02%  "

# Row 50
Set-TextValue $ws.Range("D50") "36.82"
$ws.Range("E50").Value = "  +4.93%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.514"
$ws.Range("E51").Value = "  +13.88%  "
